# Apply the "add 2022-Q3 data" edit:
#  1. Insert a new worksheet named "2022-Q3" right after "总计", before "2022-Q2".
#  2. Populate it with the fund-holding table for 2022-Q3.
#  3. Update the "总计" (summary) sheet with a new row for 2022-Q3 and renumber
#     the existing rows accordingly.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after the summary sheet.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$q3.Name = "2022-Q3"
# Match the outline settings used by the other sheets in the workbook
# (summary rows below, summary columns to the right).
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1

# Grab reference style cells from the summary sheet so the new sheet's
# header / index-column styling matches the rest of the workbook (bold,
# centered, bordered - style index 2 in the original file).
$styleHeaderSrc = $summary.Range("B1")
$styleIndexSrc = $summary.Range("A2")

# ---------------------------------------------------------------------------
# 2. Fill in the 2022-Q3 worksheet.
# ---------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

# Header row (row 1, columns B..H), styled like the other sheets' headers.
$styleHeaderSrc.Copy()
$q3.Range("B1:H1").PasteSpecial(-4122) | Out-Null
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$data = @(
    @("001508", "富国新动力灵活配置混合A", "43.83", "62.68", "3.63", "1.5910", 6),
    @("001510", "富国新动力灵活配置混合C", "12.35", "62.68", "3.63", "0.4483", 6),
    @("008372", "富国阿尔法两年持有期混合", "8.64", "72.16", "4.46", "0.3853", 5),
    @("519673", "银河康乐股票A", "1.94", "93.82", "3.62", "0.0702", 10),
    @("003284", "中邮医药健康灵活配置混合", "0.65", "76.19", "3.22", "0.0209", 7),
    @("159804", "国寿安保国证创业板中盘精选88ETF", "1.10", "98.91", "1.82", "0.0200", 9),
    @("016018", "银河康乐股票C", "0.35", "93.82", "3.62", "0.0127", 10),
    @("013920", "兴华创新医疗6个月持有混合A", "0.18", "94.83", "5.89", "0.0106", 6),
    @("013921", "兴华创新医疗6个月持有混合C", "0.05", "94.83", "5.89", "0.0029", 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rec = $data[$i]

    # Column A: numeric row index (0-based), styled like the summary sheet's
    # index column.
    $aCell = $q3.Cells.Item($row, 1)
    $styleIndexSrc.Copy()
    $aCell.PasteSpecial(-4122) | Out-Null
    $aCell.Value = $i

    # Columns B-G: text values - force text format so things like leading
    # zeros ("001508") and fixed decimal strings ("43.83") survive as text
    # instead of being coerced into numbers.
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $q3.Cells.Item($row, 2 + $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rec[$c]
    }

    # Column H: numeric rank.
    $q3.Cells.Item($row, 8).Value = $rec[6]
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert the new 2022-Q3 row and
#    renumber/shift the existing rows down by one.
# ---------------------------------------------------------------------------
# Final desired state (top to bottom):
#   row2: 0  2022-Q3   9  2.56
#   row3: 1  2022-Q2  22  7.01
#   row4: 2  2022-Q1  10  1.61
#   row5: 3  2021-Q4   8  0.71
$summaryRows = @(
    @("2022-Q3", 9, 2.56),
    @("2022-Q2", 22, 7.01),
    @("2022-Q1", 10, 1.61),
    @("2021-Q4", 8, 0.71)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $row = 2 + $i
    $rec = $summaryRows[$i]

    # Row 5 is brand new (the sheet used to only have rows 1-4), so its
    # column-A cell needs the index-column style copied over explicitly.
    if ($row -eq 5) {
        $aCell = $summary.Cells.Item($row, 1)
        $styleIndexSrc.Copy()
        $aCell.PasteSpecial(-4122) | Out-Null
    }

    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $rec[0]
    $summary.Cells.Item($row, 3).Value = $rec[1]
    $summary.Cells.Item($row, 4).Value = $rec[2]
}

$excel.CutCopyMode = $false
